$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value2 = $val
    $cell.Style = "Normal"
}

Set-TextValue 2 4 "60.306.56"
Set-TextValue 2 5 "  -2.93%  "

Set-TextValue 3 4 "3.298.79"
Set-TextValue 3 5 "  -3.67%  "

Set-TextValue 4 4 "0.999"
Set-TextValue 4 5 "  +0.01%  "

Set-TextValue 5 4 "557.19"
Set-TextValue 5 5 "  -3.94%  "

Set-TextValue 6 4 "140.89"
Set-TextValue 6 5 "  -8.95%  "

Set-TextValue 8 4 "3.297.57"
Set-TextValue 8 5 "  -3.67%  "

Set-TextValue 9 4 "0.466"
Set-TextValue 9 5 "  -3.83%  "

Set-TextValue 10 4 "7.91"
Set-TextValue 10 5 "  -1.95%  "

Set-TextValue 11 5 "  -5.42%  "

Set-TextValue 12 4 "0.406"
Set-TextValue 12 5 "  -3.10%  "

Set-TextValue 13 4 "3.865.30"
Set-TextValue 13 5 "  -3.53%  "

Set-TextValue 14 5 "  -0.17%  "

Set-TextValue 15 4 "26.56"
Set-TextValue 15 5 "  -7.54%  "

Set-TextValue 16 4 "3.307.25"
Set-TextValue 16 5 "  -3.28%  "

Set-TextValue 17 4 "0.0000163"
Set-TextValue 17 5 "  -5.46%  "

Set-TextValue 18 4 "60.261.10"
Set-TextValue 18 5 "  -3.02%  "

Set-TextValue 19 4 "6.04"
Set-TextValue 19 5 "  -7.74%  "

Set-TextValue 20 4 "13.63"
Set-TextValue 20 5 "  -5.75%  "

Set-TextValue 21 4 "8.51"
Set-TextValue 21 5 "  -5.58%  "

Set-TextValue 22 4 "373.03"
Set-TextValue 22 5 "  -2.73%  "

Set-TextValue 23 5 "  -0.04%  "

Set-TextValue 24 4 "72.19"
Set-TextValue 24 5 "  -5.12%  "

Set-TextValue 25 4 "0.530"
Set-TextValue 25 5 "  -7.48%  "

Set-TextValue 26 4 "3.437.28"
Set-TextValue 26 5 "  -3.47%  "

Set-TextValue 27 5 "  -10.08%  "

Set-TextValue 28 5 "  -2.46%  "

Set-TextValue 29 5 "  +0.30%  "

Set-TextValue 30 4 "7.01"
Set-TextValue 30 5 "  -8.85%  "

Set-TextValue 31 4 "0.999"
Set-TextValue 31 5 "  -0.05%  "

Set-TextValue 32 5 "  -4.97%  "

Set-TextValue 33 4 "7.35"
Set-TextValue 33 5 "  -6.98%  "

Set-TextValue 34 4 "22.51"
Set-TextValue 34 5 "  -3.47%  "

Set-TextValue 35 5 "  -6.26%  "

Set-TextValue 36 2 "Monero"
Set-TextValue 36 3 "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue 36 4 "165.75"
Set-TextValue 36 5 "  -1.79%  "

Set-TextValue 37 2 "NEARProtocol"
Set-TextValue 37 3 "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue 37 4 "5.01"
Set-TextValue 37 5 "  -10.32%  "

Set-TextValue 38 5 "  -5.06%  "

Set-TextValue 39 4 "6.60"
Set-TextValue 39 5 "  -5.53%  "

Set-TextValue 40 4 "3.330.56"
Set-TextValue 40 5 "  -3.73%  "

Set-TextValue 41 4 "0.0720"
Set-TextValue 41 5 "  -8.26%  "

Set-TextValue 42 4 "25.24"
Set-TextValue 42 5 "  -18.73%  "

Set-TextValue 43 4 "41.73"
Set-TextValue 43 5 "  -2.45%  "

Set-TextValue 44 4 "0.746"
Set-TextValue 44 5 "  -4.43%  "

Set-TextValue 45 4 "1.12"
Set-TextValue 45 5 "  -4.63%  "

Set-TextValue 46 4 "4.08"
Set-TextValue 46 5 "  -7.81%  "

Set-TextValue 47 4 "1.56"
Set-TextValue 47 5 "  -6.80%  "

Set-TextValue 48 5 "  -0.08%  "

Set-TextValue 49 4 "2.317.17"
Set-TextValue 49 5 "  -9.38%  "

Set-TextValue 50 2 "Cosmos"
Set-TextValue 50 3 "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue 50 4 "6.34"
Set-TextValue 50 5 "  -7.11%  "

Set-TextValue 51 2 "InjectiveProtocol"
Set-TextValue 51 3 "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue 51 4 "21.42"
Set-TextValue 51 5 "  -8.68%  "
